# Update "想去人数" (want-to-go count) values in column F
# on worksheets "展览" and "全部类型" to reflect refreshed scrape data.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitionUpdates = @{
    3  = 554
    4  = 1576
    8  = 192
    9  = 768
    12 = 364
    13 = 67
    14 = 520
    15 = 28
    16 = 6561
    17 = 30
    19 = 140
    20 = 168
    22 = 562
    23 = 15683
    24 = 1548
    25 = 19
    28 = 109
    29 = 11145
    30 = 791
    31 = 4374
    32 = 265
    35 = 310
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allTypesUpdates = @{
    3  = 554
    4  = 1576
    9  = 192
    10 = 768
    14 = 364
    15 = 67
    16 = 520
    18 = 28
    19 = 6561
    20 = 30
    22 = 140
    23 = 168
    26 = 566
    27 = 15683
    28 = 1548
    29 = 19
    32 = 109
    34 = 11145
    35 = 791
    36 = 4374
    37 = 265
    40 = 310
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
